$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Init")

$ws.Range("D5").Value = "A40"
$ws.Range("D6").Value = "B40"
$ws.Range("D7").Value = "C40"
$ws.Range("D8").Value = "G40"
$ws.Range("D9").Value = "H40"
$ws.Range("D10").Value = "I40"
$ws.Range("D11").Value = "J40"

$ws.Range("D11").Select()
